$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-26 Monday", $true, $true, $false, $false, $false, $true, 1, $false, "2026-01-27 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("344×7=", $true, $true, $false, $false, $false, $true, 1, $false, "367×9=", 2) | Out-Null
$d.Content.Find.Execute("686×6=", $true, $true, $false, $false, $false, $true, 1, $false, "422×3=", 2) | Out-Null
$d.Content.Find.Execute("986×5=", $true, $true, $false, $false, $false, $true, 1, $false, "223×3=", 2) | Out-Null
$d.Content.Find.Execute("193×2=", $true, $true, $false, $false, $false, $true, 1, $false, "800×4=", 2) | Out-Null
$d.Content.Find.Execute("996×2=", $true, $true, $false, $false, $false, $true, 1, $false, "571×5=", 2) | Out-Null
$d.Content.Find.Execute("576×2=", $true, $true, $false, $false, $false, $true, 1, $false, "183×9=", 2) | Out-Null
$d.Content.Find.Execute("405×8=", $true, $true, $false, $false, $false, $true, 1, $false, "568×4=", 2) | Out-Null
$d.Content.Find.Execute("463×8=", $true, $true, $false, $false, $false, $true, 1, $false, "666×6=", 2) | Out-Null
$d.Content.Find.Execute("998×8=", $true, $true, $false, $false, $false, $true, 1, $false, "293×9=", 2) | Out-Null
$d.Content.Find.Execute("460×5=", $true, $true, $false, $false, $false, $true, 1, $false, "944×4=", 2) | Out-Null
$d.Content.Find.Execute("934×6=", $true, $true, $false, $false, $false, $true, 1, $false, "281×5=", 2) | Out-Null
$d.Content.Find.Execute("199×7=", $true, $true, $false, $false, $false, $true, 1, $false, "988×2=", 2) | Out-Null
$d.Content.Find.Execute("613×8=", $true, $true, $false, $false, $false, $true, 1, $false, "623×7=", 2) | Out-Null
$d.Content.Find.Execute("881×8=", $true, $true, $false, $false, $false, $true, 1, $false, "903×4=", 2) | Out-Null
$d.Content.Find.Execute("667×2=", $true, $true, $false, $false, $false, $true, 1, $false, "294×4=", 2) | Out-Null
$d.Content.Find.Execute("391×3=", $true, $true, $false, $false, $false, $true, 1, $false, "927×4=", 2) | Out-Null
$d.Content.Find.Execute("238×2=", $true, $true, $false, $false, $false, $true, 1, $false, "523×4=", 2) | Out-Null
$d.Content.Find.Execute("803×8=", $true, $true, $false, $false, $false, $true, 1, $false, "840×4=", 2) | Out-Null
$d.Content.Find.Execute("333×4=", $true, $true, $false, $false, $false, $true, 1, $false, "310×4=", 2) | Out-Null
$d.Content.Find.Execute("967×2=", $true, $true, $false, $false, $false, $true, 1, $false, "415×9=", 2) | Out-Null
$d.Content.Find.Execute("131×6=", $true, $true, $false, $false, $false, $true, 1, $false, "288×2=", 2) | Out-Null
$d.Content.Find.Execute("442×7=", $true, $true, $false, $false, $false, $true, 1, $false, "467×3=", 2) | Out-Null
$d.Content.Find.Execute("379×8=", $true, $true, $false, $false, $false, $true, 1, $false, "268×4=", 2) | Out-Null
$d.Content.Find.Execute("993×2=", $true, $true, $false, $false, $false, $true, 1, $false, "265×3=", 2) | Out-Null
$d.Content.Find.Execute("208×9=", $true, $true, $false, $false, $false, $true, 1, $false, "545×9=", 2) | Out-Null
